# Daily "Updated symbol list" refresh for the cryptos sheet.
# Column D holds numeric-looking prices stored as TEXT in the source file
# (t="inlineStr"), so every Column D write is preceded by forcing the
# cell's NumberFormat to Text ("@") to keep Excel from reinterpreting the
# string as a number. Columns B/C/E are plain text and need no such care.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- price-only updates (rows 2,4-7,9) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.22"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.423"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05896"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.449"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.560"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9473"

# --- "One" jumped from rank 18 to rank 9, pushing WazirX .. CoinExToken
#     down by one row each (rows 10-18) ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01130"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1420"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07441"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03274"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03063"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09339"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.850"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001591"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04676"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- remaining scattered price-only updates ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005887"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.598"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3230"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1312"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006188"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1072"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003003"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009164"
$ws.Range("E44").Value = "43LocalTradersLCT"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005205"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7507"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002293"
